$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.808.99'
$ws.Range("E2").Value = '  +4.06%  '
$ws.Range("D3").Value = '1.913.31'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '250.61'
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '46.79'
$ws.Range("E8").Value = '  +8.26%  '
$ws.Range("E9").Value = '  +4.58%  '
$ws.Range("D10").Value = '58.14'
$ws.Range("E10").Value = '  +8.84%  '
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = '0.100'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").Value = '14.68'
$ws.Range("E13").Value = '  +8.62%  '
$ws.Range("D14").Value = '0.816'
$ws.Range("E14").Value = '  +5.07%  '
$ws.Range("D15").Value = '2.192.15'
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").Value = '5.13'
$ws.Range("E16").Value = '  +3.79%  '
$ws.Range("D17").Value = '1.907.40'
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").Value = '37.196.89'
$ws.Range("E18").Value = '  +5.05%  '
$ws.Range("D19").Value = '74.73'
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("E20").Value = '  +3.78%  '
$ws.Range("D21").Value = '13.66'
$ws.Range("E21").Value = '  +6.38%  '
$ws.Range("D22").Value = '251.61'
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  +3.36%  '
$ws.Range("D27").Value = '167.97'
$ws.Range("E27").Value = '  +1.78%  '
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("D29").Value = '18.74'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").Value = '4.63'
$ws.Range("E31").Value = '  +7.38%  '
$ws.Range("E32").Value = '  +3.95%  '
$ws.Range("D33").Value = '0.0909'
$ws.Range("E33").Value = '  +23.01%  '
$ws.Range("D34").Value = '4.35'
$ws.Range("E34").Value = '  +3.43%  '
$ws.Range("E35").Value = '  +2.51%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = '1.51'
$ws.Range("E37").Value = '  +2.38%  '
$ws.Range("D38").Value = '18.81'
$ws.Range("E38").Value = '  +56.03%  '
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("E40").Value = '  +2.72%  '
$ws.Range("D41").Value = '105.91'
$ws.Range("E41").Value = '  +8.52%  '
$ws.Range("E42").Value = '  +4.40%  '
$ws.Range("D43").Value = '17.94'
$ws.Range("E43").Value = '  +3.05%  '
$ws.Range("D44").Value = '2.89'
$ws.Range("E44").Value = '  +20.60%  '
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '1.349.90'
$ws.Range("E46").Value = '  +2.90%  '
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '0.0813'
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("D51").Value = '43.19'
$ws.Range("E51").Value = '  +1.09%  '
